# Apply the "update new orleans xlsx files" edit:
#  1. hotel_info gains a new "State" column inserted right after "Hotel_Name"
#     (and before "City"), populated with "Louisiana" for the existing row.
#  2. The sheet tab order is swapped so that "review_info" comes before
#     "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "State" column into hotel_info ----------------------
$hotelWs = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City"; shift it (and everything after) right by
# one column and put the new "State" column in its place.
$hotelWs.Columns.Item(3).Insert()

$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder the sheet tabs: review_info first, hotel_info second -------
$reviewWs = $wb.Worksheets.Item("review_info")
$reviewWs.Move($wb.Worksheets.Item(1))
